# This script applies a bugfix to the naive forecaster's YoY forecast-vector
# series: the date index advances by one period (the stale leading row is
# dropped) and every y_0_forecast / y_1_forecast value is recomputed.
#
# Net effect vs. the original sheet:
#   - Rows 2-52 are overwritten with the corrected values for the next period
#     (what used to be in rows 3-53, but with freshly computed forecasts).
#   - The now-redundant last row (53) is deleted, shrinking the used range
#     from A1:E53 to A1:E52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (row, column, new value) triples for every cell whose value changes.
$updates = @(
    @(2, 1, 39583),
    @(2, 2, 2008),
    @(2, 4, 2009),
    @(2, 5, 1.328558632615739),
    @(3, 1, 39765),
    @(3, 5, 0.4282194198276246),
    @(4, 1, 39948),
    @(4, 2, 2009),
    @(4, 3, -1.435981453719049),
    @(4, 4, 2010),
    @(4, 5, -0.7704417043119083),
    @(5, 1, 40130),
    @(5, 5, -0.8235211753995442),
    @(6, 1, 40310),
    @(6, 2, 2010),
    @(6, 3, 0.406633294022174),
    @(6, 4, 2011),
    @(6, 5, -0.1034614224434405),
    @(7, 1, 40494),
    @(7, 5, 1.121293995080253),
    @(8, 1, 40676),
    @(8, 2, 2011),
    @(8, 3, 2.185496833134781),
    @(8, 4, 2012),
    @(8, 5, 0.6652762968575532),
    @(9, 1, 40862),
    @(9, 5, 1.665250327443002),
    @(10, 1, 41044),
    @(10, 2, 2012),
    @(10, 3, 0.8574941660507873),
    @(10, 4, 2013),
    @(10, 5, 1.693469135756587),
    @(11, 1, 41228),
    @(11, 5, 1.079796209653616),
    @(12, 1, 41409),
    @(12, 2, 2013),
    @(12, 3, 0.1494732105682406),
    @(12, 4, 2014),
    @(12, 5, 0.8024032015999882),
    @(13, 1, 41592),
    @(13, 5, 1.374377011838535),
    @(14, 1, 41774),
    @(14, 2, 2014),
    @(14, 3, 1.656936590801972),
    @(14, 4, 2015),
    @(14, 5, 0.922773818606859),
    @(15, 1, 41957),
    @(15, 5, 1.310895847186577),
    @(16, 1, 42137),
    @(16, 2, 2015),
    @(16, 3, 1.346932828201242),
    @(16, 4, 2016),
    @(16, 5, 1.364302026343633),
    @(17, 1, 42321),
    @(17, 5, 1.862478303083726),
    @(18, 1, 42503),
    @(18, 2, 2016),
    @(18, 3, 1.745747589686109),
    @(18, 4, 2017),
    @(18, 5, 1.644798626926303),
    @(19, 1, 42689),
    @(19, 5, 1.639776099317536),
    @(20, 1, 42867),
    @(20, 2, 2017),
    @(20, 3, 1.843649045891893),
    @(20, 4, 2018),
    @(20, 5, 1.741128155516525),
    @(21, 1, 43053),
    @(21, 5, 2.181728312936415),
    @(22, 1, 43145),
    @(22, 2, 2018),
    @(22, 3, 2.284406789710336),
    @(22, 4, 2019),
    @(22, 5, 1.990690441067144),
    @(23, 1, 43235),
    @(23, 3, 2.463589365374652),
    @(23, 5, 2.149194501693219),
    @(24, 1, 43326),
    @(24, 3, 2.349806433215029),
    @(24, 5, 2.036910005299108),
    @(25, 1, 43418),
    @(25, 5, 2.010025322622599),
    @(26, 1, 43510),
    @(26, 2, 2019),
    @(26, 3, 1.665971362160357),
    @(26, 4, 2020),
    @(26, 5, 2.031292234149706),
    @(27, 1, 43600),
    @(27, 3, 1.332860091726285),
    @(27, 5, 1.799885362733189),
    @(28, 1, 43691),
    @(28, 3, 1.029194292875912),
    @(28, 5, 1.31420459445093),
    @(29, 1, 43783),
    @(29, 5, 0.7771393814490102),
    @(30, 1, 43875),
    @(30, 2, 2020),
    @(30, 3, 0.4126128934655471),
    @(30, 4, 2021),
    @(30, 5, 1.156986202028509),
    @(31, 1, 43966),
    @(31, 3, 0.2336391425753925),
    @(31, 5, 0.9207450904090253),
    @(32, 1, 44068),
    @(32, 3, -4.43626840667447),
    @(32, 5, -2.63419394755392),
    @(33, 1, 44159),
    @(33, 5, -2.71887004062904),
    @(34, 1, 44251),
    @(34, 2, 2021),
    @(34, 3, -2.96879819115512),
    @(34, 4, 2022),
    @(34, 5, -2.438555173006141),
    @(35, 1, 44341),
    @(35, 3, -2.010709456685855),
    @(35, 5, -1.14257141002756),
    @(36, 1, 44432),
    @(36, 3, -1.513408827666285),
    @(36, 5, 0.7106578563214505),
    @(37, 1, 44525),
    @(37, 5, 0.4582698374457683),
    @(38, 1, 44617),
    @(38, 2, 2022),
    @(38, 3, 1.154413086110817),
    @(38, 4, 2023),
    @(38, 5, -1.166698219025086),
    @(39, 1, 44706),
    @(39, 3, 1.5286818008164),
    @(39, 5, -0.7118141543333012),
    @(40, 1, 44798),
    @(40, 3, 1.618732201786743),
    @(40, 5, -0.4630595634534385),
    @(41, 1, 44890),
    @(41, 5, 1.314675624401973),
    @(42, 1, 44981),
    @(42, 2, 2023),
    @(42, 3, 0.006126408955742235),
    @(42, 4, 2024),
    @(42, 5, 0.5370151562237302),
    @(43, 1, 45071),
    @(43, 3, 0.001079933351455509),
    @(43, 5, 0.6889047703476203),
    @(44, 1, 45163),
    @(44, 3, -0.09609276733164585),
    @(44, 5, 0.5285660612534882),
    @(45, 1, 45254),
    @(45, 5, 0.1199358335146838),
    @(46, 1, 45345),
    @(46, 2, 2024),
    @(46, 3, -0.1634698065940632),
    @(46, 4, 2025),
    @(46, 5, -0.1145111565623136),
    @(47, 1, 45436),
    @(47, 3, -0.00209793826797533),
    @(47, 5, 0.286657616500996),
    @(48, 1, 45534),
    @(48, 3, -0.02761034355766023),
    @(48, 5, 0.08174908622293753),
    @(49, 1, 45618),
    @(49, 5, 0.2676745853112728),
    @(50, 1, 45713),
    @(50, 2, 2025),
    @(50, 3, 0.4725905789402463),
    @(50, 4, 2026),
    @(50, 5, 0.08564335928031852),
    @(51, 1, 45800),
    @(51, 3, 0.501314651583451),
    @(51, 5, 0.2155158706220295),
    @(52, 1, 45891),
    @(52, 3, 0.5152269879013183),
    @(52, 5, 0.3332251551730891)
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}

# Drop the obsolete trailing row (old row 53) now that its data has
# effectively moved up into row 52; this also fixes up the sheet dimension.
$ws.Rows.Item(53).Delete()
